$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.3962093333333334
$ws.Range("H2").Value = 1.188628
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.2520896666666667
$ws.Range("N2").Value = 0.7562690000000001
$ws.Range("O2").Value = 0.03491140780587004
$ws.Range("P2").Value = 0.03491140780587004
$ws.Range("Q2").Value = 0.09988027877022225
$ws.Range("R2").Value = 0.8989225089320001
$ws.Range("S2").Value = 0.03491140780587004
$ws.Range("T2").Value = 0.03491140780587004

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.3962093333333334
$ws.Range("H3").Value = 1.188628
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.312792666666667
$ws.Range("N3").Value = 3.938378
$ws.Range("O3").Value = 0.1818061039810792
$ws.Range("P3").Value = 0.1818061039810792
$ws.Range("Q3").Value = 0.5201407072648889
$ws.Range("R3").Value = 4.681266365384
$ws.Range("S3").Value = 0.1818061039810792
$ws.Range("T3").Value = 0.1818061039810792

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.3962093333333334
$ws.Range("H4").Value = 1.188628
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.655957
$ws.Range("N4").Value = 16.967871
$ws.Range("O4").Value = 0.7832824882130508
$ws.Range("P4").Value = 0.7832824882130508
$ws.Range("Q4").Value = 2.240942952332
$ws.Range("R4").Value = 20.168486570988
$ws.Range("S4").Value = 0.7832824882130508
$ws.Range("T4").Value = 0.7832824882130508
